$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Updated Ligand/Receptor-expressing cell counts (1 -> 3) and all downstream
# NATMI-computed expression/specificity metrics recalculated externally
# per commit "Natmi following Dr Hou advice".
$data = @{
    2 = @{ "E" = 3; "G" = 5.742066; "H" = 17.226198; "I" = 0.2447097919555983; "J" = 0.2447097919555983; "K" = 3; "M" = 4.765403666666667; "N" = 14.296211; "O" = 0.08284483696442199; "P" = 0.08284483696442199; "Q" = 27.363262370642; "R" = 246.269361335778; "S" = 0.02027294281815916; "T" = 0.02027294281815917 }
    3 = @{ "E" = 3; "G" = 5.742066; "H" = 17.226198; "I" = 0.2447097919555983; "J" = 0.2447097919555983; "K" = 3; "M" = 16.336489; "N" = 49.00946700000001; "O" = 0.2840040136038997; "P" = 0.2840040136038997; "Q" = 93.80519804627403; "R" = 844.2467824164661; "S" = 0.06949856308356521; "T" = 0.06949856308356521 }
    4 = @{ "E" = 3; "G" = 5.742066; "H" = 17.226198; "I" = 0.2447097919555983; "J" = 0.2447097919555983; "K" = 3; "M" = 10.06211; "N" = 30.18633; "O" = 0.174926180608571; "P" = 0.1749261806085711; "Q" = 57.77729971925999; "R" = 519.99569747334; "S" = 0.04280614926431083; "T" = 0.04280614926431085 }
    5 = @{ "E" = 3; "G" = 5.742066; "H" = 17.226198; "I" = 0.2447097919555983; "J" = 0.2447097919555983; "K" = 3; "M" = 26.358033; "N" = 79.074099; "O" = 0.4582249688231074; "P" = 0.4582249688231074; "Q" = 151.349565116178; "R" = 1362.146086045602; "S" = 0.1121321367895631; "T" = 0.1121321367895631 }
    6 = @{ "E" = 3; "G" = 6.924657666666666; "H" = 20.773973; "I" = 0.2951083350441702; "J" = 0.2951083350441703; "K" = 3; "M" = 4.765403666666667; "N" = 14.296211; "O" = 0.08284483696442199; "P" = 0.08284483696442199; "Q" = 32.99878903514477; "R" = 296.9891013163029; "S" = 0.0244482019035763; "T" = 0.02444820190357631 }
    7 = @{ "E" = 3; "G" = 6.924657666666666; "H" = 20.773973; "I" = 0.2951083350441702; "J" = 0.2951083350441703; "K" = 3; "M" = 16.336489; "N" = 49.00946700000001; "O" = 0.2840040136038997; "P" = 0.2840040136038997; "Q" = 113.1245938002657; "R" = 1018.121344202391; "S" = 0.08381195160050871; "T" = 0.08381195160050875 }
    8 = @{ "E" = 3; "G" = 6.924657666666666; "H" = 20.773973; "I" = 0.2951083350441702; "J" = 0.2951083350441703; "K" = 3; "M" = 10.06211; "N" = 30.18633; "O" = 0.174926180608571; "P" = 0.1749261806085711; "Q" = 69.67666715434332; "R" = 627.0900043890899; "S" = 0.05162217391503122; "T" = 0.05162217391503125 }
    9 = @{ "E" = 3; "G" = 6.924657666666666; "H" = 20.773973; "I" = 0.2951083350441702; "J" = 0.2951083350441703; "K" = 3; "M" = 26.358033; "N" = 79.074099; "O" = 0.4582249688231074; "P" = 0.4582249688231074; "Q" = 182.520355291703; "R" = 1642.683197625327; "S" = 0.135226007625054; "T" = 0.1352260076250541 }
    10 = @{ "E" = 3; "G" = 8.617968666666666; "H" = 25.853906; "I" = 0.3672722186578602; "J" = 0.3672722186578602; "K" = 3; "M" = 4.765403666666667; "N" = 14.296211; "O" = 0.08284483696442199; "P" = 0.08284483696442199; "Q" = 41.06809948335178; "R" = 369.612895350166; "S" = 0.03042660707627197; "T" = 0.03042660707627198 }
    11 = @{ "E" = 3; "G" = 8.617968666666666; "H" = 25.853906; "I" = 0.3672722186578602; "J" = 0.3672722186578602; "K" = 3; "M" = 16.336489; "N" = 49.00946700000001; "O" = 0.2840040136038997; "P" = 0.2840040136038997; "Q" = 140.7873503253447; "R" = 1267.086152928102; "S" = 0.1043067841840414; "T" = 0.1043067841840414 }
    12 = @{ "E" = 3; "G" = 8.617968666666666; "H" = 25.853906; "I" = 0.3672722186578602; "J" = 0.3672722186578602; "K" = 3; "M" = 10.06211; "N" = 30.18633; "O" = 0.174926180608571; "P" = 0.1749261806085711; "Q" = 86.71494870055332; "R" = 780.4345383049799; "S" = 0.06424552645345545; "T" = 0.06424552645345546 }
    13 = @{ "E" = 3; "G" = 8.617968666666666; "H" = 25.853906; "I" = 0.3672722186578602; "J" = 0.3672722186578602; "K" = 3; "M" = 26.358033; "N" = 79.074099; "O" = 0.4582249688231074; "P" = 0.4582249688231074; "Q" = 227.152702508966; "R" = 2044.374322580694; "S" = 0.1682933009440915; "T" = 0.1682933009440915 }
    14 = @{ "E" = 3; "G" = 2.180106333333333; "H" = 6.540318999999999; "I" = 0.09290965434237122; "J" = 0.09290965434237124; "K" = 3; "M" = 4.765403666666667; "N" = 14.296211; "O" = 0.08284483696442199; "P" = 0.08284483696442199; "Q" = 10.38908671458989; "R" = 93.50178043130899; "S" = 0.007697085166414546; "T" = 0.007697085166414547 }
    15 = @{ "E" = 3; "G" = 2.180106333333333; "H" = 6.540318999999999; "I" = 0.09290965434237122; "J" = 0.09290965434237124; "K" = 3; "M" = 16.336489; "N" = 49.00946700000001; "O" = 0.2840040136038997; "P" = 0.2840040136038997; "Q" = 35.61528313333034; "R" = 320.537548199973; "S" = 0.02638671473578441; "T" = 0.02638671473578442 }
    16 = @{ "E" = 3; "G" = 2.180106333333333; "H" = 6.540318999999999; "I" = 0.09290965434237122; "J" = 0.09290965434237124; "K" = 3; "M" = 10.06211; "N" = 30.18633; "O" = 0.174926180608571; "P" = 0.1749261806085711; "Q" = 21.93646973769666; "R" = 197.42822763927; "S" = 0.01625233097577353; "T" = 0.01625233097577354 }
    17 = @{ "E" = 3; "G" = 2.180106333333333; "H" = 6.540318999999999; "I" = 0.09290965434237122; "J" = 0.09290965434237124; "K" = 3; "M" = 26.358033; "N" = 79.074099; "O" = 0.4582249688231074; "P" = 0.4582249688231074; "Q" = 57.463314677509; "R" = 517.169832097581; "S" = 0.04257352346439874; "T" = 0.04257352346439874 }
}

foreach ($row in $data.Keys) {
    foreach ($col in $data[$row].Keys) {
        $ws.Range("$col$row").Value = $data[$row][$col]
    }
}
